$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numeric-looking values ("3", "4") that must stay text (to
# match the id-like "Line Number" column convention already used in the
# sheet). Format as Text first so Excel doesn't auto-convert them to
# numbers, then restore the Normal cell style so no stray formatting is
# left behind on these cells.
$ws.Range("A3:A4").NumberFormat = "@"

$ws.Range("A3").Value = "3"
$ws.Range("B3").Value = "2,priya,12346,AST,priya@tcs.com"
$ws.Range("C3").Value = "2,preya,12346,ASE-T,priya@tcs.com"

$ws.Range("A4").Value = "4"
$ws.Range("B4").Value = "3,logan,12347,ASOC,logan@tcs.com"
$ws.Range("C4").Value = "3,logan,12347,ASOC,logon@tcs.com"

$ws.Range("A3:A4").Style = "Normal"
